# Update "想去人数" (F column) counts for the matching rows on both the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.
$wb = $excel.ActiveWorkbook

# Row number -> new value for column F
$updates = @{
    2  = 316
    3  = 13911
    6  = 188
    7  = 287
    8  = 500
    10 = 91
    13 = 55
    14 = 465
    15 = 5923
    17 = 96
    18 = 987
    19 = 130
    20 = 63
    22 = 292
}

foreach ($sheetIndex in 1, 4) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
